$d = $word.ActiveDocument
$apos = [char]0x2019

function Replace-Text($old, $new) {
    $rng = $d.Content
    $found = $rng.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Text not found: $old"
    }
    $rng.Text = $new
}

# 1. Update generated timestamp
Replace-Text "Generated: 2025-09-04 10:57:32" "Generated: 2025-09-05 09:00:26"

# 2. Product overview paragraph 1
$old2 = "This product is an in-house credit card core system designed specifically for the U.S. financial services industry within the credit card sector. It serves financial institutions aiming to transition away from legacy third-party platforms to gain full ownership and control over the entire credit card lifecycle. Delivered as an integrated, real-time platform, it supports issuance of both virtual and physical credit cards, as well as comprehensive account management including credit lines, authorizations, settlements, billing, payments, interest calculations, fees, rewards, disputes, delinquency workflows, and charged-off account management. The system" + $apos + "s capability to handle end-to-end credit card processes makes it a central operational backbone for credit issuance and servicing."
$new2 = "This product is a modern, in-house credit card core platform designed specifically for the financial services industry within the United States credit card sector. It offers comprehensive credit card issuance and account management capabilities, supporting both physical and virtual card issuance in real time. The platform delivers end-to-end lifecycle management including credit line administration, product configuration, transaction authorizations, settlements, billing, payments, interest calculation, fee management, rewards, dispute handling, and delinquency workflows. It aims to fully replace legacy third-party systems with a robust, scalable infrastructure enabling full ownership and control over credit issuance and servicing processes."
Replace-Text $old2 $new2

# 3. Product overview paragraph 2 (strategic paragraph)
$old3 = "Strategically, this core system addresses the need for modernization and operational independence by embedding critical financial functions that were previously managed externally. By internalizing these capabilities, the product enables faster innovation, improved data integrity, and enhanced customer experience. It supports regulatory compliance relevant to the U.S. credit card market and positions the enterprise to adapt seamlessly to evolving business needs and competitive pressures in a dynamic financial services environment. This shift to a modern core ultimately drives greater operational efficiency and risk control across the credit portfolio. "
$new3 = "Strategically, this solution addresses the need for enhanced operational control, agility, and efficiency by transitioning from outdated legacy technologies to a modern core system. It supports sophisticated credit account management including charged-off accounts, enabling accurate lifecycle management aligned with regulatory and business requirements. This transformation allows the enterprise to streamline processes, improve customer experience, reduce dependency risks, and maintain compliance with applicable U.S. financial regulations, positioning the organization for sustainable growth and competitive differentiation in the evolving credit card market."
Replace-Text $old3 $new3

# 4. References heading (overview section) - remove trailing spaces
Replace-Text "### References  " "### References"

# 5. No external sources sentence (overview section)
Replace-Text "No external sources were used in the generation of this overview." "No external sources used."

# 6. Feature overview paragraph 1
$old6 = "The Create a Frontbook Charged-Off Account feature is a specialized subset of the broader account creation capability within the in-house credit card core system. It enables comprehensive lifecycle management of accounts that have transitioned into charged-off status due to triggers such as prolonged delinquency, customer death, or bankruptcy. This feature supports real-time processing and integration with credit line management, authorizations, settlements, billing, payments, interest calculations, fees, rewards, disputes, and delinquency workflows, ensuring seamless handling of charged-off accounts within the trade credit ecosystem."
$new6 = "The Create a Frontbook Charged-Off Account feature enables the establishment and full lifecycle management of charged-off credit card accounts within a modern in-house credit card core system. This capability is a specialized subset of the broader account creation process, activated when an account transitions from good standing to charged-off status due to triggers such as prolonged delinquency, customer death, or bankruptcy. It supports real-time processing and integration with credit issuance, account management, billing, payments, disputes, and delinquency workflows, ensuring seamless handling of charged-off accounts within the trade credit ecosystem."
Replace-Text $old6 $new6

# 7. Feature overview paragraph 2
$old7 = "This feature includes the creation, status transition, and ongoing management of charged-off accounts but excludes the initial issuance of accounts in good standing or unrelated account servicing functions. It integrates with core systems responsible for transaction posting, interest calculation, account updates, and dispute management, leveraging data inputs from daily transaction files and account status triggers. Critical constraints include compliance with financial regulations governing charged-off accounts, real-time data accuracy, and secure handling of sensitive customer information. Strategically, this feature supports the organization" + $apos + "s goal of full ownership over credit issuance and servicing by replacing legacy third-party systems with a modern, flexible core that enhances operational control, risk management, and customer lifecycle visibility. "
$new7 = "This feature includes defining and implementing the necessary functionality to create and manage charged-off accounts on the new core platform, replacing legacy third-party systems. It excludes broader account creation activities unrelated to charged-off status and focuses on the specific business rules and data flows associated with charged-off accounts. Key constraints include compliance with financial regulations, accurate interest and fee calculations, and integration with transaction posting, account updates, and dispute management systems. Strategically, this feature supports the organization" + $apos + "s goal of full ownership and control over credit card operations, improving operational efficiency, risk management, and customer servicing capabilities."
Replace-Text $old7 $new7

# 8. Final "No external sources" sentence (feature section)
Replace-Text "No external sources were used." "No external sources used."
